$d = $word.ActiveDocument

# --- Hunk 1: merge "(喀痰)クリアランスモード" + bookmark + "（振動重畳）" runs into one run ---
# Both runs share identical run formatting and are only separated by a stray
# "_GoBack" bookmark. Find/Replace across the whole (bookmark-spanning) phrase
# collapses it into a single run and drops the bookmark in the process.
$found1 = $d.Content.Find.Execute(
    "(喀痰)クリアランスモード（振動重畳）", $true, $false, $false, $false, $false,
    $true, 1, $false, "(喀痰)クリアランスモード（振動重畳）", 2)
if (-not $found1) { throw "hunk1: could not find target phrase" }

# --- Hunk 2: split "通称「火星人」" into "通称「火星人" + bookmark + "」" -----------
# Formatting is identical on both sides of the split, so simply drop a
# "_GoBack" bookmark one character before the end of the run; the engine
# keeps the surrounding run text/formatting untouched and the trailing
# character becomes its own run after the bookmark.
$r2 = $d.Content.Duplicate
$found2 = $r2.Find.Execute("通称「火星人」")
if (-not $found2) { throw "hunk2: could not find target run" }
$splitPos = $r2.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))

# --- Hunk 3: drop bold from the "↖ ... Bennett弁の特許図 ... 流体素子の概念図 ... ↗" runs ---
$r3a = $d.Content.Duplicate
$found3a = $r3a.Find.Execute("↖")
if (-not $found3a) { throw "hunk3: could not find start marker" }
$boldStart = $r3a.Start

$r3b = $d.Content.Duplicate
$found3b = $r3b.Find.Execute("↗")
if (-not $found3b) { throw "hunk3: could not find end marker" }
$boldEnd = $r3b.End

$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Font.Bold = 0

Write-Host "done"
